$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (column D) cells keep their original text formatting
# (Excel would otherwise auto-convert numeric-looking strings into numbers,
# losing formatting such as trailing zeros, thousands separators, etc.)
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.625.46'
$ws.Range('E2').Value = '  +1.84%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.548.74'
$ws.Range('E3').Value = '  +5.03%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.13'
$ws.Range('E5').Value = '  +2.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.10'
$ws.Range('E6').Value = '  +8.83%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +0.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.549.98'
$ws.Range('E9').Value = '  +5.15%  '
$ws.Range('E10').Value = '  +2.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.72'
$ws.Range('E11').Value = '  -0.93%  '
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('E13').Value = '  +3.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.39'
$ws.Range('E14').Value = '  +8.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.006.35'
$ws.Range('E15').Value = '  +5.17%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.523.21'
$ws.Range('E16').Value = '  +1.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000144'
$ws.Range('E17').Value = '  +2.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.591.99'
$ws.Range('E18').Value = '  +7.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.65'
$ws.Range('E19').Value = '  +4.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '339.88'
$ws.Range('E20').Value = '  -1.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.36'
$ws.Range('E21').Value = '  +3.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.83'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.04'
$ws.Range('E24').Value = '  +1.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.170'
$ws.Range('E25').Value = '  -1.05%  '
$ws.Range('B26').Value = 'SuiNetwork'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.55'
$ws.Range('E26').Value = '  +13.73%  '
$ws.Range('B27').Value = 'Fetch.AI'
$ws.Range('C27').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.60'
$ws.Range('E27').Value = '  +3.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.55'
$ws.Range('E28').Value = '  +4.95%  '
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.15'
$ws.Range('E30').Value = '  +10.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0829'
$ws.Range('E31').Value = '  +5.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.87'
$ws.Range('E32').Value = '  +3.27%  '
$ws.Range('E33').Value = '  +3.45%  '
$ws.Range('E34').Value = '  +11.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '422.25'
$ws.Range('E35').Value = '  +11.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.406'
$ws.Range('E36').Value = '  +2.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.11'
$ws.Range('E37').Value = '  +2.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.46'
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('B39').Value = 'USDe'
$ws.Range('C39').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.79'
$ws.Range('E40').Value = '  +5.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.63'
$ws.Range('E42').Value = '  +0.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '154.18'
$ws.Range('E43').Value = '  +6.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.79'
$ws.Range('E44').Value = '  +3.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.95'
$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.612'
$ws.Range('E46').Value = '  +3.66%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0530'
$ws.Range('E47').Value = '  +1.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0968'
$ws.Range('E48').Value = '  +1.64%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0240'
$ws.Range('E49').Value = '  +7.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.68'
$ws.Range('E50').Value = '  +4.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.82'
$ws.Range('E51').Value = '  +7.03%  '
